$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labels (row 18): "Mean increase" in D18, "Median increase" in F18 -
# both bold, matching the existing header label style used elsewhere on
# this sheet (D6/E6, D9/E9, D12/E12, D15/E15).
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true

$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

# New computed increases (row 19): percentage increase of the mean (E3)
# and median (E10) relative to fixed reference values.
$ws.Range("D19").Formula = "= ((E3 / 95.321842) * 100) - 100"
# The formula above pulls in E3's number-format (style) as a side effect
# of referencing a formatted cell; strip that back off so D19 stays
# General-formatted like the source workbook.
$ws.Range("D19").ClearFormats()

$ws.Range("F19").Formula = "= ((E10 / 95.22216) * 100) - 100"
$ws.Range("F19").ClearFormats()

# Update the active-sheet selection to match the edited area.
$ws.Range("J23").Select() | Out-Null
